$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells) ---
# A8 contains "Volume 30   Number  20" -> change trailing "20" to "22"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "22"

# C9 contains "Report Covering the Week  5/15/2023  Through  5/21/2023"
# -> dates change to 5/29/2023 and 6/4/2023
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "5/29/2023"
$c9.Characters(47, 9).Text = "6/4/2023"

# --- Crime statistics table updates (rows 14-29) ---

$ws.Range("C14").Value = "0"
$ws.Range("L14").Value = -57.142857142857
$ws.Range("M14").Value = -25
$ws.Range("N14").Value = -75
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -53.846153846153
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 76.923076923076
$ws.Range("I16").Value = 103
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = -2.830188679245
$ws.Range("L16").Value = 60.9375
$ws.Range("M16").Value = 9.574468085106
$ws.Range("N16").Value = -69.883040935672
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = 2.777777777777
$ws.Range("I17").Value = 156
$ws.Range("J17").Value = 193
$ws.Range("K17").Value = -19.170984455958
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = 45.794392523364
$ws.Range("N17").Value = -39.299610894941
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -58.333333333333
$ws.Range("I18").Value = 47
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = -41.25
$ws.Range("L18").Value = 51.612903225806
$ws.Range("M18").Value = 2.173913043478
$ws.Range("N18").Value = -78.037383177570
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -25.490196078431
$ws.Range("I19").Value = 179
$ws.Range("J19").Value = 205
$ws.Range("K19").Value = -12.682926829268
$ws.Range("L19").Value = 45.528455284552
$ws.Range("M19").Value = 94.565217391304
$ws.Range("N19").Value = 24.305555555555
$ws.Range("C20").Value = "0"
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 4
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -23.529411764705
$ws.Range("L20").Value = 8.333333333333
$ws.Range("M20").Value = 18.181818181818
$ws.Range("N20").Value = -80.152671755725
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 27
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 131
$ws.Range("H21").Value = -10.687022900763
$ws.Range("I21").Value = 520
$ws.Range("J21").Value = 626
$ws.Range("K21").Value = -16.932907348242
$ws.Range("L21").Value = 35.770234986945
$ws.Range("M21").Value = 39.410187667560
$ws.Range("N21").Value = -53.279424977538
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = -35.294117647058
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -47.619047619047
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -21.052631578947
$ws.Range("I23").Value = 68
$ws.Range("J23").Value = 71
$ws.Range("K23").Value = -4.225352112676
$ws.Range("L23").Value = 3.030303030303
$ws.Range("M23").Value = 65.853658536585
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -5.128205128205
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -8.571428571428
$ws.Range("I24").Value = 537
$ws.Range("J24").Value = 481
$ws.Range("K24").Value = 11.642411642411
$ws.Range("L24").Value = 81.418918918918
$ws.Range("M24").Value = 40.944881889763
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = -24.444444444444
$ws.Range("I25").Value = 210
$ws.Range("J25").Value = 238
$ws.Range("K25").Value = -11.764705882352
$ws.Range("L25").Value = 16.022099447513
$ws.Range("M25").Value = 12.299465240641
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 1
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 11
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 57.142857142857
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -17.241379310344
$ws.Range("L27").Value = -25
$ws.Range("G28").Value = "0"
$ws.Range("H28").Value = "***.*"
$ws.Range("I28").Value = 17
$ws.Range("K28").Value = 21.428571428571
$ws.Range("L28").Value = 21.428571428571
$ws.Range("M28").Value = -10.526315789473
$ws.Range("N28").Value = -5.555555555555
$ws.Range("G29").Value = "0"
$ws.Range("H29").Value = "***.*"
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = 23.076923076923
$ws.Range("L29").Value = 33.333333333333
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -11.111111111111